# Add exploration and test 2 rows (participants 9vshqapy, lt530p8m, yqlveap5, cbwk49s2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('9vshqapy','Training phase',1,'[''Purple'', ''Orange'', ''Orange'', ''Orange'', ''Green'', ''Purple'']','[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('9vshqapy','Training phase',2,'[''Green'', ''Green'', ''Purple'', ''Orange'', ''Purple'', ''Purple'']','[[''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('9vshqapy','Training phase',3,'[''Orange'', ''Green'', ''Purple'', ''Orange'', ''Orange'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Blue'', '''']]')
    ,@('9vshqapy','Training phase',4,'[''Orange'', ''Purple'', ''Orange'', ''Purple'', ''Green'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('9vshqapy','Training phase',5,'[''Purple'', ''Orange'', ''Green'', ''Green'', ''Orange'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('9vshqapy','Test 1',1,'[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('9vshqapy','Exploration',1,'[''Green'', ''Purple'', ''Blue'']','[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('9vshqapy','Exploration',2,'[''Green'', ''Purple'', ''Blue'']','[[''Blue'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('9vshqapy','Exploration',3,'[''Green'', ''Purple'', ''Yellow'']','[[''Red'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('9vshqapy','Exploration',4,'[''Green'', ''Purple'', ''Yellow'']','[[''Blue'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('9vshqapy','Exploration',5,'[''Green'', ''Purple'', ''Yellow'']','[[''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('lt530p8m','Training phase',1,'[''Purple'', ''Orange'', ''Orange'', ''Orange'', ''Green'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('lt530p8m','Training phase',2,'[''Green'', ''Green'', ''Purple'', ''Orange'', ''Purple'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('lt530p8m','Training phase',3,'[''Orange'', ''Green'', ''Purple'', ''Orange'', ''Orange'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('lt530p8m','Training phase',4,'[''Orange'', ''Purple'', ''Orange'', ''Purple'', ''Green'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('lt530p8m','Training phase',5,'[''Purple'', ''Orange'', ''Green'', ''Green'', ''Orange'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('lt530p8m','Test 1',1,'[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']','[[''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('lt530p8m','Exploration',1,'[''Green'', ''Purple'', ''Blue'']','[[''Red'', ''''], [''Yellow'', ''''], [''Red'', '''']]')
    ,@('yqlveap5','Training phase',1,'[''Purple'', ''Orange'', ''Orange'', ''Orange'', ''Green'', ''Purple'']','[[''Yellow'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('yqlveap5','Training phase',2,'[''Green'', ''Green'', ''Purple'', ''Orange'', ''Purple'', ''Purple'']','[[''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('yqlveap5','Training phase',3,'[''Orange'', ''Green'', ''Purple'', ''Orange'', ''Orange'', ''Green'']','[[''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('yqlveap5','Training phase',4,'[''Orange'', ''Purple'', ''Orange'', ''Purple'', ''Green'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('yqlveap5','Training phase',5,'[''Purple'', ''Orange'', ''Green'', ''Green'', ''Orange'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('yqlveap5','Test 1',1,'[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('cbwk49s2','Training phase',1,'[''Purple'', ''Orange'', ''Orange'', ''Orange'', ''Green'', ''Purple'']','[[''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Yellow'', '''']]')
    ,@('cbwk49s2','Training phase',2,'[''Green'', ''Green'', ''Purple'', ''Orange'', ''Purple'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('cbwk49s2','Training phase',3,'[''Orange'', ''Green'', ''Purple'', ''Orange'', ''Orange'', ''Green'']','[[''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Red'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Training phase',4,'[''Orange'', ''Purple'', ''Orange'', ''Purple'', ''Green'', ''Green'']','[[''Red'', ''''], [''Blue'', ''''], [''Yellow'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Training phase',5,'[''Purple'', ''Orange'', ''Green'', ''Green'', ''Orange'', ''Purple'']','[[''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Red'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Test 1',1,'[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']','[[''Red'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Red'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('cbwk49s2','Exploration phase',1,'[''Green'', ''Purple'', ''Blue'']','[[''Red'', ''''], [''Red'', ''''], [''Blue'', '''']]')
    ,@('cbwk49s2','Exploration phase',2,'[''Green'', ''Purple'', ''Blue'']','[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Exploration phase',3,'[''Green'', ''Purple'', ''Yellow'']','[[''Red'', ''''], [''Blue'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Exploration phase',4,'[''Green'', ''Purple'', ''Yellow'']','[[''Red'', ''''], [''Red'', ''''], [''Red'', '''']]')
    ,@('cbwk49s2','Exploration phase',5,'[''Green'', ''Purple'', ''Yellow'']','[[''Blue'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
    ,@('cbwk49s2','Test 2',1,'[''Green'', ''Yellow'', ''Purple'', ''Red'', ''Orange'', ''Blue'']','[[''Red'', ''''], [''Red'', ''''], [''Red'', ''''], [''Yellow'', ''''], [''Blue'', ''''], [''Blue'', '''']]')
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
}

# Widen column B so the stored OOXML width ends up at 19
# (COM's ColumnWidth setter adds a fixed 0.8333... padding when round-tripped to the XML width unit)
$ws.Range("B:B").ColumnWidth = 18.16666666666667